$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in worked hours for row 18 (B18, C18); D18 formula will recompute automatically
$ws.Range("B18").Value = 9
$ws.Range("C18").Value = 15

# Move the active selection from E17 to E19
$ws.Range("E19").Select()
